# ICTU-Kwaliteitsaanpak.pptx edit:
#   Insert new measure "M35: Het project hanteert een agile architectuuraanpak"
#   into the "Ontwikkelproces" section. This shifts every following
#   measure-slide's content down by one slide, and the slide that used to
#   be last (M33) ends up duplicated onto a brand-new trailing slide.
#
# Net effect on slide content (by slide index, 1-based):
#   20: M10  -> M35 (new text)
#   21: M16  -> M10
#   22: M28  -> M16
#   23: M30  -> M28
#   24: M34  -> M30
#   25: M27  -> M34
#   26: Organisatie -> M27
#   27: M29  -> Organisatie
#   28: M19  -> M29
#   29: M18  -> M19
#   30: M11  -> M18
#   31: M12  -> M11
#   32: M33  -> M12
#   33: (new) -> M33

$p = $ppt.ActivePresentation

function Clone-SlideContent($srcSlide, $dstSlide) {
    # Remove every shape currently on the destination slide...
    while ($dstSlide.Shapes.Count -gt 0) {
        $dstSlide.Shapes.Item(1).Delete()
    }
    # ...and replace them with faithful copies of the source slide's shapes,
    # preserving shape type (placeholder vs. textbox), geometry and every
    # paragraph/run formatting detail.
    $n = $srcSlide.Shapes.Count
    for ($i = 1; $i -le $n; $i++) {
        $srcSlide.Shapes.Item($i).Copy()
        $dstSlide.Shapes.Paste() | Out-Null
    }
}

# 1. The slide that is currently last (slide 32, "M33: ICTU organiseert
#    periodiek een gezamenlijke self-assessment...") stays the last slide of
#    the deck; duplicate it so its content survives at the new, final slide
#    position (33) once slide 32 itself is overwritten below.
$slide32 = $p.Slides.Item(32)
$slide32.Duplicate() | Out-Null

# 2. Cascade every slide's content into the next one, walking back-to-front
#    so a source slide is always read before it gets overwritten.
for ($idx = 32; $idx -ge 21; $idx--) {
    $dst = $p.Slides.Item($idx)
    $src = $p.Slides.Item($idx - 1)
    Clone-SlideContent $src $dst
}

# 3. Slide 20 gets the brand-new measure text (same shapes/structure as
#    before, only the wording changes).
$slide20 = $p.Slides.Item(20)
$slide20.Shapes.Item(1).TextFrame.TextRange.Text = "M35: Het project hanteert een agile architectuuraanpak"
$slide20.Shapes.Item(2).TextFrame.TextRange.Text = "Tijdens de voorfase verwerkt het project de door de opdrachtgever opgestelde projectstartarchitectuur (PSA) in een eerste versie van het softwarearchitectuurdocument (SAD). Tijdens de realisatiefase werkt het project het SAD bij op basis van nieuwe inzichten."
